# Published state of ETDataset on 12 February 2016
# - Refresh the dollar_per_euro exchange rate and its reference date
# - Drop the now-unused external link to coal.xlsx (rId3 / externalLink2.xml)
# - Leave the selection on the date cell that was last touched (G9)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stale external reference to coal.xlsx (this also drops the
# corresponding <externalReference> entry from workbook.xml and deletes
# xl/externalLinks/externalLink2.xml + its relationship/content-type entry).
$wb.BreakLink("coal.xlsx", 1) | Out-Null

# Updated USD/EUR running month average and its as-of date.
$ws.Range("E9").Value = 1.0965
$ws.Range("G9").Value = 42412

# Match the saved selection state (single cell G9, not the old E9:H9 block).
$ws.Range("G9").Select() | Out-Null
